$wb = $excel.ActiveWorkbook

# --- "Feedback" sheet: add a "With wind ?" column right after the
# "Starting Angle / (degrees)" column (old column G), matching the
# Starting-Angle column's formatting, and mark every experiment row "No"
# (none of the logged runs used wind). ---
$wsFeedback = $wb.Worksheets.Item("Feedback")

$wsFeedback.Columns.Item(7).Insert()
$wsFeedback.Columns.Item(7).ColumnWidth = $wsFeedback.Columns.Item(6).ColumnWidth

$wsFeedback.Range("G2").Value = "With wind ?"
$wsFeedback.Range("G4:G24").Value = "No"

# One of the starting-angle entries was filled in (row 11 was left at 0).
$wsFeedback.Range("F11").Value = 15

# --- restore/update selections and make "Feedback" the active tab ---
$wsOpenLoop = $wb.Worksheets.Item("Open Loop")
$null = $wsOpenLoop.Range("F7").Select()

$wsFeedback.Activate()
$null = $wsFeedback.Range("H16").Select()
